$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 131; everything currently at/after 131 shifts down by one.
$ws.Rows.Item(131).Insert()

# Populate the newly inserted row 131 with the new data record.
$ws.Range("A131").Value = 5
$ws.Range("B131").Value = "Macroferia Regional de Talca"
$ws.Range("C131").Value = "Maule"
$ws.Range("D131").Value = 44488
$ws.Range("E131").Value = 7
$ws.Range("F131").Value = 100112032
$ws.Range("G131").Value = "Zapallo italiano"
$ws.Range("H131").Value = "Sin especificar"
$ws.Range("I131").Value = "Primera"
$ws.Range("J131").Value = 300
$ws.Range("K131").Value = 14000
$ws.Range("L131").Value = 14000
$ws.Range("M131").Value = 14000
$ws.Range("N131").Value = "$/caja 60 unidades"
$ws.Range("O131").Value = "Región de O'Higgins"
$ws.Range("P131").Value = 233
$ws.Range("Q131").Value = 60
$ws.Range("R131").Value = "Hortaliza"
